# Applies the "(Ready for frontend)" annotations plus the bookmark /
# lastRenderedPageBreak relocation described in the commit diff.

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-MatchRange([string]$searchText) {
    # Fresh range over the whole document content each time, since earlier
    # edits shift character offsets for anything that follows them.
    $rng = $d.Range(0, $d.Content.End)
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $searchText"
    }
    return $rng
}

function Replace-WholeParagraph([string]$searchText, [string]$innerXml) {
    # Locate the text, then extend the range by one character so it also
    # covers the trailing paragraph (or cell) mark; InsertXML then swaps out
    # the complete paragraph for our replacement while leaving neighbouring
    # cells/paragraphs untouched.
    $found = Get-MatchRange $searchText
    $full = $d.Range($found.Start, $found.End + 1)
    $xml = "<w:p $wns>$innerXml</w:p>"
    $full.InsertXML($xml)
}

# 1. "Register a new user" (registerUser row): becomes bold paragraph mark,
#    plus " " and a bold "(Ready for frontend)" run.
Replace-WholeParagraph "Register a new user" (
    '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r><w:t>Register a new user</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>(Ready for frontend)</w:t></w:r>'
)

# 2. "Log in an existing user" (loginUser row): " " + bold "(Ready for
#    frontend)" run, no paragraph-mark formatting change this time.
Replace-WholeParagraph "Log in an existing user" (
    '<w:r><w:t>Log in an existing user</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>(Ready for frontend)</w:t></w:r>'
)

# 3. "List all events within startdate and enddate" (eventsByDate row).
Replace-WholeParagraph "List all events within startdate and enddate" (
    '<w:r><w:t>List all events within startdate and enddate</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>(Ready for frontend)</w:t></w:r>'
)

# 4. Drop the _GoBack bookmark from the eventsByDate table's
#    "Events array and message" cell.
Replace-WholeParagraph "Events array and message" (
    '<w:r><w:t>Events array</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
    '<w:r><w:t>message</w:t></w:r>'
)

# 5. Remove the stray <w:lastRenderedPageBreak/> from the "url/eventsByCity"
#    header cell.
Replace-WholeParagraph "url/eventsByCity" (
    '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>url</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>/</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>eventsByCity</w:t></w:r>'
)

# 6. "List all active events by city" (eventsByCity row): " " + bold
#    "(Ready for frontend)" run, followed by the relocated _GoBack bookmark.
Replace-WholeParagraph "List all active events by city" (
    '<w:r><w:t>List all active events by city</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>(Ready for frontend)</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
)
